$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.012617333333333
$ws.Range("H2").Value = 3.037852
$ws.Range("I2").Value = 0.0186050446061446
$ws.Range("J2").Value = 0.0186050446061446
$ws.Range("M2").Value = 3.916551333333333
$ws.Range("N2").Value = 11.749654
$ws.Range("O2").Value = 0.06008397860270322
$ws.Range("P2").Value = 0.06008397860270322
$ws.Range("Q2").Value = 3.96596776702311
$ws.Range("R2").Value = 35.693709903208
$ws.Range("S2").Value = 0.001117865102017931
$ws.Range("T2").Value = 0.001117865102017931
$ws.Range("G3").Value = 1.012617333333333
$ws.Range("H3").Value = 3.037852
$ws.Range("I3").Value = 0.0186050446061446
$ws.Range("J3").Value = 0.0186050446061446
$ws.Range("N3").Value = 76.19865999999999
$ws.Range("O3").Value = 0.3896556151351059
$ws.Range("P3").Value = 0.3896556151351059
$ws.Range("Q3").Value = 25.72002796425777
$ws.Range("R3").Value = 231.48025167832
$ws.Range("S3").Value = 0.007249560100623358
$ws.Range("T3").Value = 0.007249560100623358
$ws.Range("G4").Value = 1.012617333333333
$ws.Range("H4").Value = 3.037852
$ws.Range("I4").Value = 0.0186050446061446
$ws.Range("J4").Value = 0.0186050446061446
$ws.Range("M4").Value = 2.116117333333333
$ws.Range("N4").Value = 6.348352
$ws.Range("O4").Value = 0.03246344494318116
$ws.Range("P4").Value = 0.03246344494318116
$ws.Range("Q4").Value = 2.142817091100444
$ws.Range("R4").Value = 19.285353819904
$ws.Range("S4").Value = 0.0006039838412370049
$ws.Range("T4").Value = 0.0006039838412370049
$ws.Range("G5").Value = 1.012617333333333
$ws.Range("H5").Value = 3.037852
$ws.Range("I5").Value = 0.0186050446061446
$ws.Range("J5").Value = 0.0186050446061446
$ws.Range("M5").Value = 33.75239833333333
$ws.Range("N5").Value = 101.257195
$ws.Range("O5").Value = 0.5177969613190098
$ws.Range("P5").Value = 0.5177969613190098
$ws.Range("Q5").Value = 34.17826359390444
$ws.Range("R5").Value = 307.60437234514
$ws.Range("S5").Value = 0.009633635562266308
$ws.Range("T5").Value = 0.009633635562266308
$ws.Range("I6").Value = 0.7824865355506074
$ws.Range("J6").Value = 0.7824865355506075
$ws.Range("M6").Value = 3.916551333333333
$ws.Range("N6").Value = 11.749654
$ws.Range("O6").Value = 0.06008397860270322
$ws.Range("P6").Value = 0.06008397860270322
$ws.Range("Q6").Value = 166.7997279134915
$ws.Range("R6").Value = 1501.197551221424
$ws.Range("S6").Value = 0.04701490425892607
$ws.Range("T6").Value = 0.04701490425892608
$ws.Range("I7").Value = 0.7824865355506074
$ws.Range("J7").Value = 0.7824865355506075
$ws.Range("N7").Value = 76.19865999999999
$ws.Range("O7").Value = 0.3896556151351059
$ws.Range("P7").Value = 0.3896556151351059
$ws.Range("R7").Value = 9735.541301756959
$ws.Range("S7").Value = 0.3049002723449099
$ws.Range("T7").Value = 0.3049002723449099
$ws.Range("I8").Value = 0.7824865355506074
$ws.Range("J8").Value = 0.7824865355506075
$ws.Range("M8").Value = 2.116117333333333
$ws.Range("N8").Value = 6.348352
$ws.Range("O8").Value = 0.03246344494318116
$ws.Range("P8").Value = 0.03246344494318116
$ws.Range("Q8").Value = 90.12209093979021
$ws.Range("R8").Value = 811.098818458112
$ws.Range("S8").Value = 0.02540220856562771
$ws.Range("T8").Value = 0.02540220856562771
$ws.Range("I9").Value = 0.7824865355506074
$ws.Range("J9").Value = 0.7824865355506075
$ws.Range("M9").Value = 33.75239833333333
$ws.Range("N9").Value = 101.257195
$ws.Range("O9").Value = 0.5177969613190098
$ws.Range("P9").Value = 0.5177969613190098
$ws.Range("Q9").Value = 1437.461271224102
$ws.Range("R9").Value = 12937.15144101692
$ws.Range("S9").Value = 0.4051691503811438
$ws.Range("T9").Value = 0.4051691503811439
$ws.Range("G10").Value = 10.82599466666667
$ws.Range("H10").Value = 32.477984
$ws.Range("I10").Value = 0.198908419843248
$ws.Range("J10").Value = 0.198908419843248
$ws.Range("M10").Value = 3.916551333333333
$ws.Range("N10").Value = 11.749654
$ws.Range("O10").Value = 0.06008397860270322
$ws.Range("P10").Value = 0.06008397860270322
$ws.Range("Q10").Value = 42.40056384639288
$ws.Range("R10").Value = 381.605074617536
$ws.Range("S10").Value = 0.01195120924175922
$ws.Range("T10").Value = 0.01195120924175922
$ws.Range("G11").Value = 10.82599466666667
$ws.Range("H11").Value = 32.477984
$ws.Range("I11").Value = 0.198908419843248
$ws.Range("J11").Value = 0.198908419843248
$ws.Range("N11").Value = 76.19865999999999
$ws.Range("O11").Value = 0.3896556151351059
$ws.Range("P11").Value = 0.3896556151351059
$ws.Range("Q11").Value = 274.9754289223822
$ws.Range("R11").Value = 2474.77886030144
$ws.Range("S11").Value = 0.0775057826895727
$ws.Range("T11").Value = 0.0775057826895727
$ws.Range("G12").Value = 10.82599466666667
$ws.Range("H12").Value = 32.477984
$ws.Range("I12").Value = 0.198908419843248
$ws.Range("J12").Value = 0.198908419843248
$ws.Range("M12").Value = 2.116117333333333
$ws.Range("N12").Value = 6.348352
$ws.Range("O12").Value = 0.03246344494318116
$ws.Range("P12").Value = 0.03246344494318116
$ws.Range("Q12").Value = 22.90907496470756
$ws.Range("R12").Value = 206.181674682368
$ws.Range("S12").Value = 0.006457252536316445
$ws.Range("T12").Value = 0.006457252536316445
$ws.Range("G13").Value = 10.82599466666667
$ws.Range("H13").Value = 32.477984
$ws.Range("I13").Value = 0.198908419843248
$ws.Range("J13").Value = 0.198908419843248
$ws.Range("M13").Value = 33.75239833333333
$ws.Range("N13").Value = 101.257195
$ws.Range("O13").Value = 0.5177969613190098
$ws.Range("P13").Value = 0.5177969613190098
$ws.Range("Q13").Value = 365.4032843438756
$ws.Range("R13").Value = 3288.62955909488
$ws.Range("S13").Value = 0.1029941753755996
$ws.Range("T13").Value = 0.1029941753755996
